# "Generate Report for Archive"
# The two localization files (322627cf-...md and 6ce311f6-...md) swapped
# their reporting order (6ce311f6 now listed before 322627cf on every
# sheet), and the 6ce311f6 file's status flipped from "Ready for handoff"
# to "In Translation" (322627cf stays "Ready for handoff").
#
# The hyperlink targets (Address/URLs) stay anchored to the same cell
# refs they always were; only the cell text / hyperlink display text
# changes - matching the source diff, where relationship ids keep
# pointing at the same targets but the <v> shared-string index (and the
# hyperlink `display=`) attached to each ref is what moves.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet --------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("A2").Value = "6ce311f6-437e-467c-a86a-aacdd1524fc0.md"
$ws.Range("B2").Value = "In Translation"
$ws.Range("C2").Value = "In Translation"
$ws.Range("A3").Value = "322627cf-f28b-4ad4-bcf4-d45a3baf76c2.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') { $h.TextToDisplay = "6ce311f6-437e-467c-a86a-aacdd1524fc0.md" }
    elseif ($addr -eq '$A$3') { $h.TextToDisplay = "322627cf-f28b-4ad4-bcf4-d45a3baf76c2.md" }
}

# ---- zh-cn sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("A2").Value = "6ce311f6-437e-467c-a86a-aacdd1524fc0.md"
$ws.Range("B2").Value = "In Translation"
$ws.Range("C2").Value = "6ce311f6-437e-467c-a86a-aacdd1524fc0.978dfe207d0544e8408b73645b7a47a5bab8eeb6.zh-cn.xlf"
$ws.Range("D2").Value = "2016-03-10 07:50:59"

$ws.Range("A3").Value = "322627cf-f28b-4ad4-bcf4-d45a3baf76c2.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "322627cf-f28b-4ad4-bcf4-d45a3baf76c2.2f1186afabb1ffe9f55f429b19fad9d79f07a323.zh-cn.xlf"
$ws.Range("D3").Value = "2016-03-10 07:50:18"

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') { $h.TextToDisplay = "6ce311f6-437e-467c-a86a-aacdd1524fc0.md" }
    elseif ($addr -eq '$C$2') { $h.TextToDisplay = "6ce311f6-437e-467c-a86a-aacdd1524fc0.978dfe207d0544e8408b73645b7a47a5bab8eeb6.zh-cn.xlf" }
    elseif ($addr -eq '$A$3') { $h.TextToDisplay = "322627cf-f28b-4ad4-bcf4-d45a3baf76c2.md" }
    elseif ($addr -eq '$C$3') { $h.TextToDisplay = "322627cf-f28b-4ad4-bcf4-d45a3baf76c2.2f1186afabb1ffe9f55f429b19fad9d79f07a323.zh-cn.xlf" }
}

# ---- de-de sheet -------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("A2").Value = "6ce311f6-437e-467c-a86a-aacdd1524fc0.md"
$ws.Range("B2").Value = "In Translation"
$ws.Range("C2").Value = "6ce311f6-437e-467c-a86a-aacdd1524fc0.978dfe207d0544e8408b73645b7a47a5bab8eeb6.de-de.xlf"
$ws.Range("D2").Value = "2016-03-10 07:51:08"

$ws.Range("A3").Value = "322627cf-f28b-4ad4-bcf4-d45a3baf76c2.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "322627cf-f28b-4ad4-bcf4-d45a3baf76c2.2f1186afabb1ffe9f55f429b19fad9d79f07a323.de-de.xlf"
$ws.Range("D3").Value = "2016-03-10 07:50:27"

foreach ($h in $ws.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') { $h.TextToDisplay = "6ce311f6-437e-467c-a86a-aacdd1524fc0.md" }
    elseif ($addr -eq '$C$2') { $h.TextToDisplay = "6ce311f6-437e-467c-a86a-aacdd1524fc0.978dfe207d0544e8408b73645b7a47a5bab8eeb6.de-de.xlf" }
    elseif ($addr -eq '$A$3') { $h.TextToDisplay = "322627cf-f28b-4ad4-bcf4-d45a3baf76c2.md" }
    elseif ($addr -eq '$C$3') { $h.TextToDisplay = "322627cf-f28b-4ad4-bcf4-d45a3baf76c2.2f1186afabb1ffe9f55f429b19fad9d79f07a323.de-de.xlf" }
}
